$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1525.8889
$ws.Range("J19").Value = 178
$ws.Range("L19").Value = 178
$ws.Range("N19").Value = -528
$ws.Range("H33").Value = 193
$ws.Range("I33").Value = 193
$ws.Range("K33").Value = 193
$ws.Range("M33").Value = 36
$ws.Range("H93").Value = 10556
$ws.Range("J93").Value = 10556
$ws.Range("L93").Value = 10556
$ws.Range("N93").Value = -15548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5102.25
$ws.Range("I45").Value = 5102.25
$ws.Range("K45").Value = 5102.25
$ws.Range("M45").Value = -4725.25
$ws.Range("H103").Value = 25000
$ws.Range("J103").Value = 25000
$ws.Range("L103").Value = 25000
$ws.Range("N103").Value = -27344

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2025
$ws.Range("I94").Value = 2025
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2025
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1574
$ws.Range("N94").Value = $null
$ws.Range("H100").Value = 18606
$ws.Range("J100").Value = 18606
$ws.Range("L100").Value = 18606
$ws.Range("N100").Value = -20770
$ws.Range("H134").Value = 8242.200000000001
$ws.Range("I134").Value = 2606
$ws.Range("K134").Value = 7818
$ws.Range("M134").Value = -5283

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 325.6
$ws.Range("I7").Value = 396.75
$ws.Range("K7").Value = 396.75
$ws.Range("M7").Value = -283.75
$ws.Range("H16").Value = 10000
$ws.Range("I16").Value = 10000
$ws.Range("K16").Value = 10000
$ws.Range("M16").Value = -9713
$ws.Range("H44").Value = 5000
$ws.Range("I44").Value = 5000
$ws.Range("K44").Value = 5000
$ws.Range("M44").Value = -4558
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32080
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 10000
$ws.Range("K113").Value = 10000
$ws.Range("M113").Value = -7830
$ws.Range("H134").Value = 9600
$ws.Range("J134").Value = 12333.333
$ws.Range("L134").Value = 36999.999
$ws.Range("N134").Value = -42069.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = $null
$ws.Range("H39").Value = 5000
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
$ws.Range("H45").Value = 600
$ws.Range("J45").Value = 600
$ws.Range("L45").Value = 1800
$ws.Range("N45").Value = -2864
$ws.Range("H55").Value = 4500
$ws.Range("J55").Value = 4500
$ws.Range("L55").Value = 13500
$ws.Range("N55").Value = -13854
$ws.Range("H114").Value = 894.1667
$ws.Range("I114").Value = 973
$ws.Range("J114").Value = 500
$ws.Range("K114").Value = 2919
$ws.Range("L114").Value = 1500
$ws.Range("M114").Value = 335
$ws.Range("N114").Value = -8008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 3600
$ws.Range("I27").Value = 5000
$ws.Range("K27").Value = 5000
$ws.Range("M27").Value = -4834
$ws.Range("H31").Value = 1031
$ws.Range("I31").Value = 1031
$ws.Range("K31").Value = 1031
$ws.Range("M31").Value = -739
$ws.Range("H37").Value = 1031
$ws.Range("I37").Value = 1031
$ws.Range("K37").Value = 1031
$ws.Range("M37").Value = -754
$ws.Range("H97").Value = 927.7143
$ws.Range("I97").Value = 755
$ws.Range("J97").Value = 996.8
$ws.Range("K97").Value = 755
$ws.Range("L97").Value = 996.8
$ws.Range("M97").Value = -259
$ws.Range("N97").Value = -1988.8
$ws.Range("H113").Value = 2344.875
$ws.Range("I113").Value = 2394.1428
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2394.1428
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -224.1428000000001
$ws.Range("N113").Value = -6340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 14500
$ws.Range("I4").Value = 14500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 14500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -14387
$ws.Range("N4").Value = $null
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = $null
$ws.Range("H26").Value = 1200
$ws.Range("J26").Value = 1200
$ws.Range("L26").Value = 1200
$ws.Range("N26").Value = -1790
$ws.Range("H28").Value = 14500
$ws.Range("I28").Value = 14500
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 14500
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -14268
$ws.Range("N28").Value = $null
$ws.Range("H37").Value = 14500
$ws.Range("I37").Value = 14500
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 14500
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -14393
$ws.Range("N37").Value = $null
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = $null
$ws.Range("H134").Value = 90000
$ws.Range("J134").Value = 90000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -100140
$ws.Range("H136").Value = 10714.143
$ws.Range("J136").Value = 22000
$ws.Range("L136").Value = 66000
$ws.Range("N136").Value = -71100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 2998
$ws.Range("I7").Value = 2998
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2998
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2885
$ws.Range("N7").Value = $null
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = $null
$ws.Range("H26").Value = 5000
$ws.Range("J26").Value = 5000
$ws.Range("L26").Value = 5000
$ws.Range("N26").Value = -5586
$ws.Range("H58").Value = 22100
$ws.Range("I58").Value = 22100
$ws.Range("K58").Value = 22100
$ws.Range("M58").Value = -21792
$ws.Range("H68").Value = 36439
$ws.Range("J68").Value = 36439
$ws.Range("L68").Value = 36439
$ws.Range("N68").Value = -38061
$ws.Range("H71").Value = 36439
$ws.Range("J71").Value = 36439
$ws.Range("L71").Value = 109317
$ws.Range("N71").Value = -117429
$ws.Range("H101").Value = 27500
$ws.Range("J101").Value = 27500
$ws.Range("L101").Value = 27500
$ws.Range("N101").Value = -33990
$ws.Range("H103").Value = 47067.332
$ws.Range("J103").Value = 47067.332
$ws.Range("L103").Value = 47067.332
$ws.Range("N103").Value = -49411.332
$ws.Range("H122").Value = 2000.75
$ws.Range("I122").Value = 2000.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6002.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3552.25
$ws.Range("N122").Value = $null
$ws.Range("H126").Value = 1249.5
$ws.Range("I126").Value = 499
$ws.Range("K126").Value = 1497
$ws.Range("M126").Value = 973

